# SCD0011 .. SCD0016 batch rename pass for this workbook (SCD0015-001).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new test-case id.
$ws.Name = "SCD0015"

# Update the TC_ID cell from the old Jira-style id to the new SCD id.
$ws.Range("B2").Value = "SCD0015-001"

# Best-fit the TC_ID column now that its content is wider.
$ws.Columns.Item(2).ColumnWidth = 11.6

# Move the active selection/view back to the top of the sheet.
$ws.Range("B3").Select()
